$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in cell A1
$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 10:44"

# Updated per-country statistics (B=Casos totales, C=Nuevos casos,
# D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)

# Row 6 - Rusia
$ws.Range("B6").Value = 545458
$ws.Range("C6").Value = 8248
$ws.Range("D6").Value = 294306
$ws.Range("E6").Value = 243868
$ws.Range("G6").Value = 193
$ws.Range("H6").Value = 7284

# Row 7 - India
$ws.Range("B7").Value = 344407
$ws.Range("C7").Value = 1381
$ws.Range("D7").Value = 180460
$ws.Range("E7").Value = 154026
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 9921

# Row 13 - Alemania
$ws.Range("D13").Value = 173100
$ws.Range("E13").Value = 6059

# Row 33 - Singapur
$ws.Range("B33").Value = 40969
$ws.Range("C33").Value = 151
$ws.Range("E33").Value = 10577

# Row 40 - Polonia
$ws.Range("B40").Value = 30195
$ws.Range("C40").Value = 407
$ws.Range("D40").Value = 14654
$ws.Range("E40").Value = 14269
$ws.Range("G40").Value = 16
$ws.Range("H40").Value = 1272

# Row 44 - Oman
$ws.Range("B44").Value = 25269
$ws.Range("C44").Value = 745
$ws.Range("D44").Value = 11089
$ws.Range("E44").Value = 14066
$ws.Range("G44").Value = 6
$ws.Range("H44").Value = 114

# Row 54 - Austria
$ws.Range("B54").Value = 17189
$ws.Range("C54").Value = 54
$ws.Range("D54").Value = 16089
$ws.Range("E54").Value = 419
$ws.Range("G54").Value = 3
$ws.Range("H54").Value = 681

# Row 61 - Moldavia
$ws.Range("D61").Value = 6901
$ws.Range("E61").Value = 4561
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 417

# Row 65 - Chequia
$ws.Range("B65").Value = 10066
$ws.Range("C65").Value = 2
$ws.Range("D65").Value = 7300
$ws.Range("E65").Value = 2436

# Row 87 - El Salvador
$ws.Range("B87").Value = 3941
$ws.Range("C87").Value = 115
$ws.Range("D87").Value = 2041
$ws.Range("E87").Value = 1824

# Row 102 - Estonia
$ws.Range("B102").Value = 1975
$ws.Range("C102").Value = 1
$ws.Range("D102").Value = 1728
$ws.Range("E102").Value = 178

# Row 111 - Eslovaquia
$ws.Range("D111").Value = 1426
$ws.Range("E111").Value = 98

# Row 146 - Togo
$ws.Range("B146").Value = 532
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 308
$ws.Range("E146").Value = 211

# Row 149 - Estado de Palestina
$ws.Range("B149").Value = 506
$ws.Range("C149").Value = 1
$ws.Range("E149").Value = 88

# Row 173 - Brunei
$ws.Range("E173").Value = 0
$ws.Range("G173").Value = 1
$ws.Range("H173").Value = 3
